# Add new "MgCa Coretop modelled temperature" column (AP), shifting the
# existing "MgCa Temperature anomaly_Original - Coretop" (old AP) and
# "MgCa Temperature anomaly_BAYMAG - Coretop" (old AQ) columns one to the
# right (-> AQ, AR). Excel's column Insert naturally carries over header
# styling (bold + centered) and all formula/reference shifts, and expands
# the sheet's used-range/dimension for us.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("AP").Insert()

# Populate the new header label for the inserted column.
$ws.Range("AP1").Value = "MgCa Coretop modelled temperature"
$ws.Range("AD2").Value = 1.45
$ws.Range("AE2").Value = 2.046592968040042
$ws.Range("AF2").Value = 1.705640728040043
$ws.Range("AD3").Value = 6.78
$ws.Range("AE3").Value = 7.059325807359484
$ws.Range("AF3").Value = 6.400870697359484
$ws.Range("AD4").Value = 10.67
$ws.Range("AE4").Value = 6.299139543321397
$ws.Range("AF4").Value = 5.942828103321398
$ws.Range("AD5").Value = 17.76
$ws.Range("AE5").Value = 3.835014970567492
$ws.Range("AF5").Value = 3.092009860567494
$ws.Range("AG5").Value = 4.476264970567492
$ws.Range("AH5").Value = 9.226684970567494
$ws.Range("AP5").Value = 18.6015
$ws.Range("AQ5").Value = 3.638499999999997
$ws.Range("AR5").Value = 8.388919999999999
$ws.Range("AD6").Value = 17.76
$ws.Range("AE6").Value = 2.811264970567493
$ws.Range("AF6").Value = 2.568582730567492
$ws.Range("AQ6").ClearContents()
$ws.Range("AD7").Value = 22.22
$ws.Range("AG7").Value = -1.791507281679124
$ws.Range("AH7").Value = 2.938683627411777
$ws.Range("AP7").Value = 24.5493
$ws.Range("AQ7").Value = -4.12202727
$ws.Range("AR7").Value = 0.6081636400000008
$ws.Range("AD8").Value = 27.59
$ws.Range("AG8").Value = -3.07169423421227
$ws.Range("AH8").Value = -3.54476090087887
$ws.Range("AP8").Value = 27.5334
$ws.Range("AQ8").Value = -3.020066670000002
$ws.Range("AR8").Value = -3.493133329999999
$ws.Range("AD9").Value = 25.31
$ws.Range("AE9").Value = 1.428633694118922
$ws.Range("AF9").Value = 2.690959834118921
$ws.Range("AD10").Value = 18.45
$ws.Range("AE10").Value = 8.433762016464843
$ws.Range("AF10").Value = 8.766535096464846
$ws.Range("AD11").Value = 19.22
$ws.Range("AE11").Value = 6.568540089925129
$ws.Range("AF11").Value = 6.355465859925129
$ws.Range("AD12").Value = 13.35
$ws.Range("AE12").Value = 4.399195263974608
$ws.Range("AF12").Value = 3.70739198397461
$ws.Range("AG12").Value = -1.68524917602539
$ws.Range("AH12").Value = -0.8866664487526901
$ws.Range("AP12").Value = 13.8771
$ws.Range("AR12").Value = -1.418517270000001
$ws.Range("Y13").Value = "van der Weijst and Peterse (Unpublished data)"
$ws.Range("AD13").Value = 27.31
$ws.Range("AG13").Value = 0.2875847710503443
$ws.Range("AQ13").ClearContents()
$ws.Range("AR13").ClearContents()
$ws.Range("Y14").Value = "van der Weijst and Peterse (Unpublished data)"
$ws.Range("AD14").Value = 25.3
$ws.Range("AE14").Value = 2.621301235622827
$ws.Range("AF14").Value = 3.538465645622825
$ws.Range("W15").Value = 24.789288925
$ws.Range("AD15").Value = 20.74
$ws.Range("AG15").Value = 4.049772560796445
$ws.Range("AH15").Value = 5.313683635796444
$ws.Range("AP15").Value = 22.0667
$ws.Range("AR15").Value = 3.986499999999999
$ws.Range("AD16").Value = 16.09
$ws.Range("AE16").Value = 8.211365169949005
$ws.Range("AF16").Value = 7.638039979949003
$ws.Range("AD17").Value = 16.51
$ws.Range("AE17").Value = 9.45864292054253
$ws.Range("AF17").Value = 9.328026210542532
$ws.Range("AD18").Value = 16.25
$ws.Range("AE18").Value = 8.049847157796222
$ws.Range("AF18").Value = 7.492214147796222
$ws.Range("AQ18").ClearContents()
$ws.Range("AD19").Value = 16.93
$ws.Range("AE19").Value = 3.513918789302302
$ws.Range("AF19").Value = 2.650603959302302
$ws.Range("AD20").Value = 28.44
$ws.Range("AE20").Value = 0.4449956936306414
$ws.Range("AF20").Value = 2.11188818363064
$ws.Range("AG20").Value = -0.7341709730360613
$ws.Range("AH20").Value = -1.308487639702662
$ws.Range("AP20").Value = 27.6266
$ws.Range("AQ20").Value = 0.07923333000000099
$ws.Range("AR20").Value = -0.4950833299999999
$ws.Range("AD21").Value = 27.73
$ws.Range("AG21").Value = -3.698869323730506
$ws.Range("AH21").Value = -0.9651359903971048
$ws.Range("AP21").Value = 25.1871
$ws.Range("AQ21").Value = -1.15376667
$ws.Range("AR21").Value = 1.579966669999997
$ws.Range("AD22").Value = 7.21
$ws.Range("AE22").Value = 2.867240799096137
$ws.Range("AF22").Value = 2.281823642096137
$ws.Range("W23").Value = 29.37598672
$ws.Range("X23").Value = 30.8326333333333
$ws.Range("AD23").Value = 29.09
$ws.Range("AG23").Value = 0.282595301542969
$ws.Range("AH23").Value = 1.739241914876267
$ws.Range("AP23").Value = 28.7843
$ws.Range("AQ23").Value = 0.5916867199999984
$ws.Range("AR23").Value = 2.048333333333296
$ws.Range("AD24").Value = 23.78
$ws.Range("AE24").Value = 1.95582495462023
$ws.Range("AF24").Value = 1.72504086462023
$ws.Range("AD25").Value = 26.36
$ws.Range("AE25").Value = 1.644249301486543
$ws.Range("AF25").Value = 3.298275921486542
$ws.Range("AD26").Value = 27.01
$ws.Range("AG26").Value = -0.01264004177517464
$ws.Range("AH26").Value = -1.062940041775175
$ws.Range("AP26").Value = 27.2019
$ws.Range("AQ26").Value = -0.2018999999999984
$ws.Range("AR26").Value = -1.252199999999998
$ws.Range("AD27").Value = 27.68
$ws.Range("AG27").Value = -1.182914225260415
$ws.Range("AH27").Value = -0.8977142252604153
$ws.Range("AP27").Value = 28.2441
$ws.Range("AD28").Value = 26.1
$ws.Range("AE28").Value = 1.813480360243055
$ws.Range("AF28").Value = 2.742746040243055
$ws.Range("AD29").Value = 25.31
$ws.Range("AG29").Value = -0.8899967108832421
$ws.Range("AH29").Value = 0.4471032891167575
$ws.Range("AP29").Value = 25.7765
$ws.Range("AQ29").Value = -1.356499999999997
$ws.Range("AR29").Value = -0.01939999999999742
$ws.Range("AD30").Value = 15.49
$ws.Range("AE30").Value = 0.1586896209733055
$ws.Range("AF30").Value = -0.5350416490266952
$ws.Range("AD31").Value = 10.66
$ws.Range("AE31").Value = 1.921620837741427
$ws.Range("AF31").Value = 1.289436677741428
$ws.Range("AD32").Value = 9.73
$ws.Range("AE32").Value = 4.129912928701174
$ws.Range("AF32").Value = 3.736367028701174
$ws.Range("AD33").Value = 14.54
$ws.Range("AE33").Value = 5.109476860894096
$ws.Range("AF33").Value = 4.369063940894096
